$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (shifts old D -> E, old E -> F)
$ws.Range("D1").EntireColumn.Insert()

# Set the header for the new column D, and fix up the shifted "score (%)" header to "score"
$ws.Range("D1").Value = "match_query"
$ws.Range("E1").Value = "score"

# Populate the new column D with concatenated first name + last name + email
$data = @(
    @(2, "April", "Gonzalez", "user8@example.com"),
    @(3, "Colton", "Collins", "user9@example.com"),
    @(4, "Joseph", "Yang", "user10@example.com"),
    @(5, "Michelle", "Logan", "user11@example.com"),
    @(6, "Mary", "Pollard", "user12@example.com"),
    @(7, "Gail", "Harris", "user13@example.com"),
    @(8, "Lauren", "Fletcher", "user14@example.com"),
    @(9, "Alan", "Liu", "user15@example.com"),
    @(10, "Douglas", "Williams", "user16@example.com"),
    @(11, "Daniel", "Smith", "user17@example.com")
)

foreach ($row in $data) {
    $r = $row[0]
    $first = $row[1]
    $last = $row[2]
    $email = $row[3]
    $ws.Cells.Item($r, 4).Value = "$first$last$email"
}
